$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header C1: audioFalse -> currentPhase
$ws.Range("C1").Value = "currentPhase"

# Column C (was per-row audio file paths) now both rows reference "train2P2"
$ws.Range("C2").Value = "train2P2"
$ws.Range("C3").Value = "train2P2"
